$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'black cigarette pants',
    'black cigarette pants women',
    'black clothes for women',
    'black clothing tape',
    'black compression',
    'black compression hose for women',
    'black compression leggings for women',
    'black compression leggings for women high waist',
    'black compression pants',
    'black compression short',
    'black compression shorts',
    'black compression shorts for women',
    'black compression shorts women',
    'black compression tight',
    'black compression tights',
    'black compression tights men',
    'black compression tights women',
    'black conpression pants',
    'black control top leggings',
    'black cross leggings',
    'black culture art',
    'black cut leggings',
    'black define band',
    'black extra long leggings',
    'black fashion tights',
    'black fatigue pants for men',
    'black football leggings',
    'black frisbee',
    'black full length leggings',
    'black gym pants',
    'black high rise leggings',
    'black high rise leggings for women',
    'black high waist pants for women',
    'black high waist tights',
    'black high waist tights for women',
    'black high waisted soft leggings for women',
    'black high waisted tights',
    'black high wasted leggings',
    'black holy leggings',
    'black hose for women',
    'black insulated skirt',
    'black kinesiology tape',
    'black knee high hose',
    'black knee high stocking',
    'black knee high stockings plus size',
    'black knee high tights',
    'black knee highs for women',
    'black knee stockings',
    'black leggibgs women',
    'black legging',
    'black legging capris for women',
    'black legging long',
    'black leggings 3 pack',
    'black leggings capri for women',
    'black leggings capri plus size',
    'black leggings control top',
    'black leggings design',
    'black leggings for women athletic',
    'black leggings for women capri',
    'black leggings for women compression',
    'black leggings for women control top',
    'black leggings for women high waist',
    'black leggings for women long',
    'black leggings for women long length',
    'black leggings gym',
    'black leggings ladies',
    'black leggings long',
    'black leggings long length',
    'black leggings lululemon',
    'black leggings running',
    'black leggings short',
    'black leggings short length',
    'black leggings tall',
    'black leggings tall and curvy',
    'black leggings victoria secret',
    'black leggings with control top',
    'black leggings with rips',
    'black leggings women',
    'black leggings women short',
    'black legings',
    'black line',
    'black long leggings for women',
    'black long pants',
    'black lululemon leggings',
    'black lululemon shorts',
    'black martial arts pants',
    'black maternity tights',
    'black medium',
    'black mild cigarette',
    'black mountain band',
    'black n',
    'black neoprene',
    'black neoprene fabric',
    'black nike pros women',
    'black nylon tights',
    'black opaque tights women',
    'black panel leggings',
    'black pants',
    'black pants high waist',
    'black pants red stripe women'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
